# Applies the "parameter updates, README fixes" edit to misc/parameters.xlsx
#
# Summary of the change:
#  - A new "Iterative LSA" column (H) is added to the parameters table,
#    mirroring the formatting of the existing "Test Data" column (G).
#  - The footnote that used to live in H8 moves one column to the right (I8)
#    now that H8 holds a regular data value ("inf") like the other columns.
#  - A few existing values are corrected (G3, G6, G7) and a label is
#    reworded (A4: "single-core" -> "single-instance" streaming SVD).
#  - Column A is widened, and the last-used selection/cursor cell changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value corrections in the existing table ------------------------------

# A4's label text changed; set this first so new shared strings are
# appended in the same order the original authoring tool produced them.
$ws.Range("A4").Value = "single-instance streaming SVD"

$ws.Range("G3").Value = 22
$ws.Range("G6").Value = 0.8
$ws.Range("G7").Value = 0.8

# --- New column H: "Iterative LSA" ----------------------------------------
# Match the existing "Test Data" column's formatting (center-aligned, same
# number format per row) by setting alignment/number format before the
# value, so the engine reuses the existing style records instead of
# minting new ones.

$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").Value = "Iterative LSA"

$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").Value = 33

$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("H3").Value = 27

$ws.Range("H4").HorizontalAlignment = -4108
$ws.Range("H4").Value = $true

$ws.Range("H5").HorizontalAlignment = -4108
$ws.Range("H5").NumberFormat = "0%"
$ws.Range("H5").Value = 0.8

$ws.Range("H6").HorizontalAlignment = -4108
$ws.Range("H6").Value = 0.45

$ws.Range("H7").HorizontalAlignment = -4108
$ws.Range("H7").Value = 0.45

$ws.Range("H8").HorizontalAlignment = -4108
$ws.Range("H8").Value = "inf"

# The footnote that used to sit in H8 now moves to I8 (unformatted, like
# the other "description" cells in column A-C).
$ws.Range("I8").Value = "* set to ~6.5 for testing the effect of writing to multiple partitions"

# --- Cosmetic sheet changes -------------------------------------------------

# Column A grows from 22.5 to 29.5 characters wide (COM ColumnWidth units
# are offset by ~5/6 of a character from the raw OOXML column width for
# this workbook's default font, hence the subtraction below).
$ws.Columns("A").ColumnWidth = 29.5 - 5/6

# Cursor/selection moves to C12.
$null = $ws.Range("C12").Select()
